# Apply scheduled-runner profit recalculations to Jenova_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Cells.Item(40, 8).Value = 7430
$ws.Cells.Item(40, 9).Value = 6380.6
$ws.Cells.Item(40, 10).Value = 7804.7856
$ws.Cells.Item(40, 11).Value = 6380.6
$ws.Cells.Item(40, 12).Value = 7804.7856
$ws.Cells.Item(40, 13).Value = -6205.6
$ws.Cells.Item(40, 14).Value = -8154.7856

# Row 41
$ws.Cells.Item(41, 9).Value = 601.6
$ws.Cells.Item(41, 11).Value = 601.6
$ws.Cells.Item(41, 13).Value = -161.6

# Row 98
$ws.Cells.Item(98, 8).Value = 2212.077
$ws.Cells.Item(98, 9).Value = 1575.697
$ws.Cells.Item(98, 10).Value = 5712.1665
$ws.Cells.Item(98, 11).Value = 1575.697
$ws.Cells.Item(98, 12).Value = 5712.1665
$ws.Cells.Item(98, 13).Value = -77.69699999999989
$ws.Cells.Item(98, 14).Value = -8708.166499999999

# Row 122
$ws.Cells.Item(122, 8).Value = 2212.077
$ws.Cells.Item(122, 9).Value = 1575.697
$ws.Cells.Item(122, 10).Value = 5712.1665
$ws.Cells.Item(122, 11).Value = 4727.090999999999
$ws.Cells.Item(122, 12).Value = 17136.4995
$ws.Cells.Item(122, 13).Value = -2277.090999999999
$ws.Cells.Item(122, 14).Value = -22036.4995

# Row 132
$ws.Cells.Item(132, 8).Value = 2306.375
$ws.Cells.Item(132, 9).Value = 2073.366
$ws.Cells.Item(132, 11).Value = 6220.098
$ws.Cells.Item(132, 13).Value = -3690.098

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 1112.2
$ws.Cells.Item(2, 9).Value = 1206
$ws.Cells.Item(2, 10).Value = 737
$ws.Cells.Item(2, 11).Value = 1206
$ws.Cells.Item(2, 12).Value = 737
$ws.Cells.Item(2, 13).Value = -1093
$ws.Cells.Item(2, 14).Value = -963

# Row 32
$ws.Cells.Item(32, 8).Value = 4594.7886
$ws.Cells.Item(32, 9).Value = 4594.7886
$ws.Cells.Item(32, 11).Value = 4594.7886
$ws.Cells.Item(32, 13).Value = -4307.7886

# Row 97
$ws.Cells.Item(97, 8).Value = 857.4761999999999
$ws.Cells.Item(97, 9).Value = 757.94116
$ws.Cells.Item(97, 10).Value = 1280.5
$ws.Cells.Item(97, 11).Value = 757.94116
$ws.Cells.Item(97, 12).Value = 1280.5
$ws.Cells.Item(97, 13).Value = -261.94116
$ws.Cells.Item(97, 14).Value = -2272.5

# Row 116
$ws.Cells.Item(116, 8).Value = 1112.2
$ws.Cells.Item(116, 9).Value = 1206
$ws.Cells.Item(116, 10).Value = 737
$ws.Cells.Item(116, 11).Value = 1206
$ws.Cells.Item(116, 12).Value = 737
$ws.Cells.Item(116, 13).Value = 1088
$ws.Cells.Item(116, 14).Value = -5325

# Row 122
$ws.Cells.Item(122, 8).Value = 3267.7073
$ws.Cells.Item(122, 9).Value = 1954.8462
$ws.Cells.Item(122, 11).Value = 5864.5386
$ws.Cells.Item(122, 13).Value = -3414.5386

# Row 132
$ws.Cells.Item(132, 8).Value = 2032.3572
$ws.Cells.Item(132, 9).Value = 1448.919
$ws.Cells.Item(132, 11).Value = 4346.757000000001
$ws.Cells.Item(132, 13).Value = -1816.757000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 1112.2
$ws.Cells.Item(3, 9).Value = 1206
$ws.Cells.Item(3, 10).Value = 737
$ws.Cells.Item(3, 11).Value = 1206
$ws.Cells.Item(3, 12).Value = 737
$ws.Cells.Item(3, 13).Value = -1092
$ws.Cells.Item(3, 14).Value = -965

# Row 134
$ws.Cells.Item(134, 8).Value = 29176.736
$ws.Cells.Item(134, 9).Value = 1756.75
$ws.Cells.Item(134, 10).Value = 175416.67
$ws.Cells.Item(134, 11).Value = 5270.25
$ws.Cells.Item(134, 12).Value = 526250.01
$ws.Cells.Item(134, 13).Value = -2735.25
$ws.Cells.Item(134, 14).Value = -531320.01

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Cells.Item(99, 8).Value = 5283.5454
$ws.Cells.Item(99, 9).Value = 3751.2727
$ws.Cells.Item(99, 11).Value = 3751.2727
$ws.Cells.Item(99, 13).Value = -2253.2727

# Row 126
$ws.Cells.Item(126, 8).Value = 5283.5454
$ws.Cells.Item(126, 9).Value = 3751.2727
$ws.Cells.Item(126, 11).Value = 11253.8181
$ws.Cells.Item(126, 13).Value = -8783.8181

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Cells.Item(13, 8).Value = 2141.6667
$ws.Cells.Item(13, 10).Value = 4225
$ws.Cells.Item(13, 12).Value = 12675
$ws.Cells.Item(13, 14).Value = -13011

# Row 114
$ws.Cells.Item(114, 8).Value = 2160
$ws.Cells.Item(114, 9).Value = 404.75
$ws.Cells.Item(114, 10).Value = 3330.1667
$ws.Cells.Item(114, 11).Value = 1214.25
$ws.Cells.Item(114, 12).Value = 9990.500100000001
$ws.Cells.Item(114, 13).Value = 2039.75
$ws.Cells.Item(114, 14).Value = -16498.5001

# Row 129
$ws.Cells.Item(129, 8).Value = 17545354
$ws.Cells.Item(129, 10).Value = 2134.9167
$ws.Cells.Item(129, 12).Value = 6404.750100000001
$ws.Cells.Item(129, 14).Value = -16404.7501

# Row 131
$ws.Cells.Item(131, 8).Value = 2816.1167
$ws.Cells.Item(131, 9).Value = 974.5
$ws.Cells.Item(131, 10).Value = 3020.7407
$ws.Cells.Item(131, 11).Value = 2923.5
$ws.Cells.Item(131, 12).Value = 9062.222099999999
$ws.Cells.Item(131, 13).Value = 2116.5
$ws.Cells.Item(131, 14).Value = -19142.2221

# Row 137
$ws.Cells.Item(137, 8).Value = 2558.9583
$ws.Cells.Item(137, 9).Value = 1544.6842
$ws.Cells.Item(137, 11).Value = 4634.0526
$ws.Cells.Item(137, 13).Value = 465.9474

# Row 138
$ws.Cells.Item(138, 8).Value = 11113259
$ws.Cells.Item(138, 9).Value = 2189.8572
$ws.Cells.Item(138, 11).Value = 6569.571599999999
$ws.Cells.Item(138, 13).Value = -1429.571599999999

$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Cells.Item(95, 8).Value = 47737.5
$ws.Cells.Item(95, 10).Value = 47737.5
$ws.Cells.Item(95, 12).Value = 47737.5
$ws.Cells.Item(95, 14).Value = -53229.5

# Row 102
$ws.Cells.Item(102, 8).Value = 2083.2273
$ws.Cells.Item(102, 9).Value = 497
$ws.Cells.Item(102, 11).Value = 497
$ws.Cells.Item(102, 13).Value = 1125

# Row 113
$ws.Cells.Item(113, 8).Value = 597753.8
$ws.Cells.Item(113, 9).Value = 1114232.2
$ws.Cells.Item(113, 11).Value = 1114232.2
$ws.Cells.Item(113, 13).Value = -1112062.2

# Row 122
$ws.Cells.Item(122, 8).Value = 3959.3809
$ws.Cells.Item(122, 9).Value = 1649.4
$ws.Cells.Item(122, 11).Value = 4948.200000000001
$ws.Cells.Item(122, 13).Value = -2498.200000000001

# Row 132
$ws.Cells.Item(132, 8).Value = 337409.16
$ws.Cells.Item(132, 9).Value = 479594.72
$ws.Cells.Item(132, 11).Value = 1438784.16
$ws.Cells.Item(132, 13).Value = -1436254.16

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 3337223.5
$ws.Cells.Item(40, 9).Value = 5003398
$ws.Cells.Item(40, 11).Value = 5003398
$ws.Cells.Item(40, 13).Value = -5003262

# Row 55
$ws.Cells.Item(55, 8).Value = 704.75
$ws.Cells.Item(55, 9).Value = 217.07692
$ws.Cells.Item(55, 11).Value = 217.07692
$ws.Cells.Item(55, 13).Value = -44.07692

# Row 61
$ws.Cells.Item(61, 8).Value = 4539.7085
$ws.Cells.Item(61, 9).Value = 3998.3333
$ws.Cells.Item(61, 11).Value = 3998.3333
$ws.Cells.Item(61, 13).Value = -3796.3333

# Row 82
$ws.Cells.Item(82, 8).Value = 489.22223
$ws.Cells.Item(82, 9).Value = 500.42856
$ws.Cells.Item(82, 10).Value = 450
$ws.Cells.Item(82, 11).Value = 500.42856
$ws.Cells.Item(82, 12).Value = 450
$ws.Cells.Item(82, 13).Value = -139.42856
$ws.Cells.Item(82, 14).Value = -1172

# Row 85
$ws.Cells.Item(85, 8).Value = 489.22223
$ws.Cells.Item(85, 9).Value = 500.42856
$ws.Cells.Item(85, 10).Value = 450
$ws.Cells.Item(85, 11).Value = 500.42856
$ws.Cells.Item(85, 12).Value = 450
$ws.Cells.Item(85, 13).Value = 747.5714399999999
$ws.Cells.Item(85, 14).Value = -2946

# Row 113
$ws.Cells.Item(113, 8).Value = 4539.7085
$ws.Cells.Item(113, 9).Value = 3998.3333
$ws.Cells.Item(113, 11).Value = 3998.3333
$ws.Cells.Item(113, 13).Value = -1828.3333

$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Cells.Item(74, 8).Value = 6824.364
$ws.Cells.Item(74, 9).Value = 4657
$ws.Cells.Item(74, 10).Value = 7637.125
$ws.Cells.Item(74, 11).Value = 4657
$ws.Cells.Item(74, 12).Value = 7637.125
$ws.Cells.Item(74, 13).Value = -3721
$ws.Cells.Item(74, 14).Value = -9509.125

# Row 77
$ws.Cells.Item(77, 8).Value = 6824.364
$ws.Cells.Item(77, 9).Value = 4657
$ws.Cells.Item(77, 10).Value = 7637.125
$ws.Cells.Item(77, 11).Value = 13971
$ws.Cells.Item(77, 12).Value = 22911.375
$ws.Cells.Item(77, 13).Value = -9291
$ws.Cells.Item(77, 14).Value = -32271.375

# Row 81
$ws.Cells.Item(81, 8).Value = 1637.619
$ws.Cells.Item(81, 9).Value = 1155.7333
$ws.Cells.Item(81, 10).Value = 2842.3333
$ws.Cells.Item(81, 11).Value = 2311.4666
$ws.Cells.Item(81, 12).Value = 5684.6666
$ws.Cells.Item(81, 13).Value = -1250.4666
$ws.Cells.Item(81, 14).Value = -7806.6666

# Row 84
$ws.Cells.Item(84, 8).Value = 1637.619
$ws.Cells.Item(84, 9).Value = 1155.7333
$ws.Cells.Item(84, 10).Value = 2842.3333
$ws.Cells.Item(84, 11).Value = 11557.333
$ws.Cells.Item(84, 12).Value = 28423.333
$ws.Cells.Item(84, 13).Value = -6253.333000000001
$ws.Cells.Item(84, 14).Value = -39031.333

# Row 123
$ws.Cells.Item(123, 8).Value = 88950
$ws.Cells.Item(123, 10).Value = 88950
$ws.Cells.Item(123, 12).Value = 88950
$ws.Cells.Item(123, 14).Value = -98750

